$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "21R01A67E6"
$ws.Range("B14").Value = "21r01a67e6"
$ws.Range("C14").Value = "CMRIT25_21R01A67E6"
$ws.Range("D14").Value = "21r01a67e6"
$ws.Range("E14").Value = "r_21r01a67e6"
$ws.Range("F14").Value = "21R01A67E6"

$ws.Range("A14:F14").Select()
